$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.764.06"
$ws.Range("E2").Value = "  -2.15%  "
$ws.Range("D3").Value = "1.875.68"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  -1.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.687"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.26%  "
$ws.Range("E7").Value = "  -1.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.48"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.347"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "50.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0740"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0968"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("D14").Value = "2.146.02"
$ws.Range("E14").Value = "  -2.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.714"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").Value = "1.871.18"
$ws.Range("E17").Value = "  -2.29%  "
$ws.Range("D18").Value = "34.740.22"
$ws.Range("E18").Value = "  -2.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").Value = "0.0₃0823"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "247.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.01%  "
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("E25").Value = "  +3.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.99%  "
$ws.Range("E30").Value = "  -3.81%  "
$ws.Range("D31").Value = "4.128.36"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +14.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0580"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("E37").Value = "  -5.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.835"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.22%  "
$ws.Range("E39").Value = "  -3.18%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.55%  "
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0210"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("E44").Value = "  -4.76%  "
$ws.Range("D45").Value = "1.292.64"
$ws.Range("E45").Value = "  -4.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.27%  "
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("E48").Value = "  -2.26%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0764"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.94%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("B51").Value = "Gas"
$ws.Range("C51").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "12.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.11%  "
